# Update the "想去人数" (want-to-go count) values in column F
# for rows 2, 3, 7, 8, 10, 15, 19, 22 on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 295
    3  = 13818
    7  = 273
    8  = 492
    10 = 86
    15 = 5771
    19 = 88
    22 = 230
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
